$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (Beta) values
$ws.Range("C2").Value = 16.27395536621195
$ws.Range("E2").Value = 0.006275575909959944
$ws.Range("F2").Value = 5.419143735856542
$ws.Range("G2").Value = 3.489393161372961
$ws.Range("H2").Value = 7.358364726262872
$ws.Range("I2").Value = 0.03122526930932872
$ws.Range("J2").Value = 0.00650207751422318
$ws.Range("K2").Value = 0.0725564117308637
$ws.Range("L2").Value = 0.006852645963075444
$ws.Range("M2").Value = 0.002890481616023623
$ws.Range("N2").Value = 0.01262880953653676

# Update existing row 3 (Gamma) values
$ws.Range("C3").Value = 0.3589057182506037
$ws.Range("D3").Value = 0.3038166771491592
$ws.Range("E3").Value = 0.3557873748505794
$ws.Range("F3").Value = 0.501976381459472
$ws.Range("G3").Value = 0.01311796782735086
$ws.Range("H3").Value = 1.077735895611366
$ws.Range("I3").Value = 0.4636924753817018
$ws.Range("J3").Value = 0.0121250613178008
$ws.Range("K3").Value = 0.9946740427033118
$ws.Range("L3").Value = 0.5209492973534806
$ws.Range("M3").Value = 0.01377028691417379
$ws.Range("N3").Value = 1.11414161745914

# Add new row 4 (Beta + Gamma)
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 16.63286108446255
$ws.Range("D4").Value = 0.3068349094013313
$ws.Range("E4").Value = 0.3620629507605393
$ws.Range("F4").Value = 5.921120117316013
$ws.Range("G4").Value = 3.502511129200312
$ws.Range("H4").Value = 8.436100621874237
$ws.Range("I4").Value = 0.4949177446910304
$ws.Range("J4").Value = 0.01862713883202399
$ws.Range("K4").Value = 1.067230454434176
$ws.Range("L4").Value = 0.5278019433165559
$ws.Range("M4").Value = 0.01666076853019741
$ws.Range("N4").Value = 1.126770426995677
